# Intelligibility_Score.xlsx - fix IndexOutOfBound error by re-syncing
# the SENTENCES / INPUT_SENTENCE / NAME / INTELLIGIBILITY_SCORE table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("B1").Value = "SENTENCES"
$ws.Range("C1").Value = "INPUT_SENTENCE"
$ws.Range("D1").Value = "NAME"
$ws.Range("E1").Value = "INTELLIGIBILITY_SCORE"

# --- Data rows ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "I think I'm getting better."
$ws.Range("C2").Value = "Test1"
$ws.Range("D2").Value = "P1_W2_S1"
$ws.Range("E2").Value = 0.125

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "You want him to do well"
$ws.Range("C3").Value = "test2"
$ws.Range("D3").Value = "P1_W2_S2"
$ws.Range("E3").Value = 0.14285714285714279

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Big muscles are not necessarily strong ones"
$ws.Range("C4").Value = "test3"
$ws.Range("D4").Value = "P1_W2_S3"
$ws.Range("E4").Value = 0.125

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "he is capable and willing to make decisions."
$ws.Range("C5").Value = "test4"
$ws.Range("D5").Value = "P1_W2_S4"
$ws.Range("E5").Value = 0.1224489795918367

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "We picked grapes for wine"
$ws.Range("C6").Value = "test5"
$ws.Range("D6").Value = "P1_W1_S1"
$ws.Range("E6").Value = 0.1333333333333333

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "The ballet is about to begin."
$ws.Range("C7").Value = "test6"
$ws.Range("D7").Value = "P1_W1_S2"
$ws.Range("E7").Value = 0.1764705882352941

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "You're used to being on the field."
$ws.Range("C8").Value = "test7"
$ws.Range("D8").Value = "P1_W1_S3"
$ws.Range("E8").Value = 0.1538461538461539

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Enjoy the fair weather while in the tropics."
$ws.Range("C9").Value = "test8"
$ws.Range("D9").Value = "P1_W1_S4"
$ws.Range("E9").Value = 0.081632653061224483
